$pres = $ppt.ActivePresentation

function Update-LinkText {
    param($SlideIndex, $ShapeIndex, $OldUrl, $NewUrl)

    $slide = $pres.Slides.Item($SlideIndex)
    $shape = $slide.Shapes.Item($ShapeIndex)
    $tr = $shape.TextFrame.TextRange
    $idx = $tr.Text.IndexOf($OldUrl)
    $sub = $tr.Characters($idx + 1, $OldUrl.Length)
    $sub.Text = $NewUrl
}

# Slide 6 - Content Placeholder 2 - unsafe-code language specification link
Update-LinkText 6 2 "https://docs.microsoft.com/en-us/dotnet/csharp/language-reference/language-specification/unsafe-code" "https://learn.microsoft.com/en-us/dotnet/csharp/language-reference/language-specification/unsafe-code"

# Slide 24 - Content Placeholder 2 - data type ranges + windows data types links
Update-LinkText 24 2 "https://docs.microsoft.com/en-us/cpp/cpp/data-type-ranges" "https://learn.microsoft.com/en-us/cpp/cpp/data-type-ranges"
Update-LinkText 24 2 "https://docs.microsoft.com/en-us/windows/win32/winprog/windows-data-types" "https://learn.microsoft.com/en-us/windows/win32/winprog/windows-data-types"

# Slide 29 - Content Placeholder 2 - safehandle, criticalfinalizerobject, constrained-execution-regions links
Update-LinkText 29 2 "https://docs.microsoft.com/en-us/dotnet/api/system.runtime.interopservices.safehandle" "https://learn.microsoft.com/en-us/dotnet/api/system.runtime.interopservices.safehandle"
Update-LinkText 29 2 "https://docs.microsoft.com/en-us/dotnet/api/system.runtime.constrainedexecution.criticalfinalizerobject" "https://learn.microsoft.com/en-us/dotnet/api/system.runtime.constrainedexecution.criticalfinalizerobject"
Update-LinkText 29 2 "https://docs.microsoft.com/en-us/dotnet/framework/performance/constrained-execution-regions" "https://learn.microsoft.com/en-us/dotnet/framework/performance/constrained-execution-regions"

Write-Output "Done updating links"
